# Applies the cryptos list update for commit "Updated cryptos list on Wed Sep 25 06:17:02 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a bare number-looking string must be forced to
# Text format first, otherwise COM auto-converts them to numeric cells and
# mangles formatting (e.g. "0.110" -> 0.11, "1.00" -> 1). This matches how
# the source report keeps the Price column as literal text.

$ws.Range('D2').Value = '64.096.11'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = '2.623.03'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.64'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.41'
$ws.Range('E6').Value = '  +2.86%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.588'
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.109'
$ws.Range('E9').Value = '  +1.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.67'
$ws.Range('E10').Value = '  +1.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.384'
$ws.Range('E11').Value = '  +6.17%  '
$ws.Range('E12').Value = '  -0.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.64'
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('D14').Value = '3.094.54'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').Value = '63.964.39'
$ws.Range('E15').Value = '  +1.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000149'
$ws.Range('E16').Value = '  +3.19%  '
$ws.Range('D17').Value = '2.629.38'
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.20'
$ws.Range('E18').Value = '  +7.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.65'
$ws.Range('E19').Value = '  +3.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '350.39'
$ws.Range('E20').Value = '  +2.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.92'
$ws.Range('E21').Value = '  +0.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.72'
$ws.Range('E23').Value = '  +2.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.41'
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.74'
$ws.Range('E25').Value = '  +14.98%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.25'
$ws.Range('E26').Value = '  +6.66%  '
$ws.Range('B27').Value = 'Fetch.AI'
$ws.Range('C27').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.69'
$ws.Range('E27').Value = '  +2.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.19'
$ws.Range('E28').Value = '  +4.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.164'
$ws.Range('E29').Value = '  +0.99%  '
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '539.29'
$ws.Range('E31').Value = '  -1.39%  '
$ws.Range('E32').Value = '  +1.68%  '
$ws.Range('D33').Value = '0.0₃0852'
$ws.Range('E33').Value = '  +6.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.76'
$ws.Range('E34').Value = '  +0.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.29'
$ws.Range('E35').Value = '  +0.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '168.17'
$ws.Range('E36').Value = '  +1.29%  '
$ws.Range('E37').Value = '  +1.22%  '
$ws.Range('B38').Value = 'FirstDigitalUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.97'
$ws.Range('E39').Value = '  +5.30%  '
$ws.Range('E40').Value = '  +2.82%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '168.56'
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.02'
$ws.Range('E43').Value = '  +0.92%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.92'
$ws.Range('E44').Value = '  +4.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0589'
$ws.Range('E45').Value = '  +3.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.54'
$ws.Range('E46').Value = '  -3.69%  '
$ws.Range('E47').Value = '  +0.95%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.00'
$ws.Range('E48').Value = '  +13.89%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0246'
$ws.Range('E49').Value = '  +1.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0967'
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.23'
$ws.Range('E51').Value = '  +3.23%  '
